$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename a handful of labels in column A to match the restructured
# "csv2sql" scraps directory (Cluster -> DC, BizDev -> Manager, IPSGO -> plain).
$ws.Range("A1").Value = "DCName"
$ws.Range("A2").Value = "DCID"
$ws.Range("A12").Value = "ManagerSelectDate"
$ws.Range("A29").Value = "Select"
$ws.Range("A30").Value = "ManagerSelect"
$ws.Range("A32").Value = "Engineer"
$ws.Range("A33").Value = "ManagerRepresentative"

# Put the active selection on A2, matching the resaved workbook view.
$ws.Range("A2").Select()
